$d = $word.ActiveDocument

# --- Locate the hidden "**ID__..." marker run in the first paragraph ---
$idOld = "**ID__AFFARS_5309_topic_14__ID**"
$idNew = "**ID__AFFARS_5309_405_2__ID**"

$findRange = $d.Content
$found = $findRange.Find.Execute($idOld, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$idStart = $findRange.Start
$idEnd = $findRange.End

# The paragraph that owns the marker run (needed for its pPr/border/indent).
$p1 = $d.Range($idStart, $idEnd).Paragraphs(1)

# Add a paragraph border (top/left/bottom/right), each with 5pt "space"
# padding, i.e. <w:pBdr><w:top w:space="5"/><w:left w:space="5"/>
# <w:bottom w:space="5"/><w:right w:space="5"/></w:pBdr>
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Change the paragraph's left indent from 120 twips (6pt) to 225 twips
# (11.25pt) -> <w:ind w:left="225"/>
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# The marker text is immediately followed, in its own run, by a single
# trailing space (the rest of the paragraph). Delete that separate run's
# character entirely so the paragraph ends up with one run holding just
# the new ID text (no trailing space run left behind).
$trailingRange = $d.Range($idEnd, $idEnd + 1)
if ($trailingRange.Text -eq " ") {
    $trailingRange.Delete()
}

# Finally, swap in the new ID text.
$idRange = $d.Range($idStart, $idEnd)
$idRange.Text = $idNew
